# aggiornamento fino a 27/05
# Append new daily rows (256-269, dates 2021-05-14 .. 2021-05-27) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row, date-serial (col A), col B (nuovi pos.), col C (somma mobile 7gg.), col D (somma mobile 7gg. per 100mila abitanti)
$rows = @(
    @(256, 44330, 3, 4, 101.7293997965412),
    @(257, 44331, 3, 7, 178.0264496439471),
    @(258, 44332, 0, 7, 178.0264496439471),
    @(259, 44333, 1, 7, 178.0264496439471),
    @(260, 44334, 0, 7, 178.0264496439471),
    @(261, 44335, 0, 7, 178.0264496439471),
    @(262, 44336, 0, 7, 178.0264496439471),
    @(263, 44337, 0, 4, 101.7293997965412),
    @(264, 44338, 0, 1, 25.4323499491353),
    @(265, 44339, 0, 1, 25.4323499491353),
    @(266, 44340, 0, 0, 0),
    @(267, 44341, 0, 0, 0),
    @(268, 44342, 0, 0, 0),
    @(269, 44343, 0, 0, 0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Column A carries the date-style formatting (s="2") used throughout the
# sheet; copy that formatting from the last pre-existing row (255) down
# across the newly added date cells so the style index is reused rather
# than a new one minted.
$ws.Range("A255").Copy()
$ws.Range("A256:A269").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Appended rows 256-269 through 2021-05-27"
